$d = $word.ActiveDocument
$wrNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:r='http://schemas.openxmlformats.org/officeDocument/2006/relationships'"

function Set-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $xml = "<w:p " + $wrNs + ">" + $innerXml + "</w:p>"
    $p.Range.InsertXML($xml)
}

function Restyle-Hyperlink($paraIndex, $displayText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $f = $p.Range
    $f.Find.Execute($displayText, $true) | Out-Null
    $f.Style = "Hyperlink"
}

# --- Q2: "Q2: Simple Directmedia Layer" -> split out "Directmedia" as a flagged misspelling
$inner = "<w:r><w:t xml:space='preserve'>Q2: Simple </w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellStart'/>"
$inner = $inner + "<w:r><w:t>Directmedia</w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellEnd'/>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> Layer</w:t></w:r>"
Set-ParagraphXml 4 $inner

# --- Q3: "Q3: Released under the zlib license <link>" -> split out "zlib"; keep hyperlink intact
$inner = "<w:r><w:t xml:space='preserve'>Q3: Released under the </w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellStart'/>"
$inner = $inner + "<w:r><w:t>zlib</w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellEnd'/>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> license </w:t></w:r>"
$inner = $inner + "<w:hyperlink r:id='rId6' w:history='1'><w:r><w:t>https://www.zlib.net/zlib_license.html</w:t></w:r></w:hyperlink>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> </w:t></w:r>"
Set-ParagraphXml 5 $inner
Restyle-Hyperlink 5 "https://www.zlib.net/zlib_license.html"

# --- Q4: "Q4: Supports windows, mac, linux, ios and android" -> split out "linux" and "ios"
$inner = "<w:r><w:t xml:space='preserve'>Q4: Supports windows, mac, </w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellStart'/>"
$inner = $inner + "<w:r><w:t>linux</w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellEnd'/>"
$inner = $inner + "<w:r><w:t xml:space='preserve'>, </w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellStart'/>"
$inner = $inner + "<w:r><w:t>ios</w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellEnd'/>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> and android</w:t></w:r>"
Set-ParagraphXml 6 $inner

# --- Source Code bullet: split out "wanna" (preserve the ListParagraph/numbering pPr)
$listPPr = "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr>"
$inner = $listPPr
$inner = $inner + "<w:r><w:t xml:space='preserve'>Source Code: The source code of the for people who </w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellStart'/>"
$inner = $inner + "<w:r><w:t>wanna</w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellEnd'/>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> work on the library itself</w:t></w:r>"
Set-ParagraphXml 10 $inner

# --- Dev Libraries bullet: split out "devs"; keep the two trailing runs untouched (preserve pPr)
$inner = $listPPr
$inner = $inner + "<w:r><w:t xml:space='preserve'>Dev Libraries: Code for game </w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellStart'/>"
$inner = $inner + "<w:r><w:t>devs</w:t></w:r>"
$inner = $inner + "<w:proofErr w:type='spellEnd'/>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> to use for</w:t></w:r>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> making</w:t></w:r>"
$inner = $inner + "<w:r><w:t xml:space='preserve'> their games</w:t></w:r>"
Set-ParagraphXml 12 $inner

# --- Append the two new paragraphs at the end of the document (new note + hyperlink to key codes)
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.Text = "Used this link to get key codes"
$endRange.InsertParagraphAfter()

$endRange = $d.Range($d.Content.End, $d.Content.End)
$url = "http://wiki.libsdl.org/SDL_Keycode?highlight=%28%5CbCategoryEnum%5Cb%29%7C%28SDLEnumTemplate%29"
$d.Hyperlinks.Add($endRange, $url) | Out-Null

# Hyperlinks.Add leaves a stray empty trailing run in the paragraph; rebuild the
# paragraph cleanly (same r:id) so it matches a normal hand-authored hyperlink run.
$lastParaIdx = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastParaIdx)
$oxml = $p.Range.WordOpenXML
$oxml -match 'hyperlink r:id="(rId\d+)"' | Out-Null
$newRid = $matches[1]

$inner = "<w:hyperlink r:id='" + $newRid + "' w:history='1'><w:r><w:t>" + $url + "</w:t></w:r></w:hyperlink>"
Set-ParagraphXml $lastParaIdx $inner
Restyle-Hyperlink $lastParaIdx $url

Write-Host "done"
